$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by one day
# (45177 -> 45178) for every data row (rows 2 through 200).
$range = $ws.Range("C2:C200")
foreach ($cell in $range.Cells) {
    if ($cell.Value2 -eq 45177) {
        $cell.Value2 = 45178
    }
}
